# Applies the fk_datasheet.docx edit described by the diff:
#  1. Merge the three "Laser Position on " / "Workbench (" / "Measured) "
#     runs into a single run with the same combined text (no content
#     change, just a run-merge).
#  2. Convert the four joint-angle test-point arrays from degrees to the
#     radians values used in the updated table.

$d = $word.ActiveDocument

# 1) Header cell: merge the three runs into one (identical concatenated
#    text, including the trailing non-breaking space).
$oldHeader = "Laser Position on Workbench (Measured)" + [char]0x00A0
$newHeader = "Laser Position on Workbench (Measured)" + [char]0x00A0
$found0 = $d.Content.Find.Execute($oldHeader, $false, $false, $false, $false, $false, $true, 1, $false, $newHeader, 2)
Write-Output "header merge found=$found0"

# 2) Test point arrays: degrees -> radians
$found1 = $d.Content.Find.Execute("[0, -45, 0, 45, -90, 60]", $false, $false, $false, $false, $false, $true, 1, $false, "[0, -0.758, 0, 0.758, -1.571, 1.048]", 2)
Write-Output "row1 found=$found1"

$found2 = $d.Content.Find.Execute("[-30, -60, 80, -10, -90, -30]", $false, $false, $false, $false, $false, $true, 1, $false, "[-0.524, -1.048, 1.396, -0.175, -1.571,        -0.524]", 2)
Write-Output "row2 found=$found2"

$found3 = $d.Content.Find.Execute("[30 -70 80 -10 -90 10]", $false, $false, $false, $false, $false, $true, 1, $false, "[0.524, -1.222, 1.396, -0.175, -1.571, 0.175]", 2)
Write-Output "row3 found=$found3"

$found4 = $d.Content.Find.Execute("[-30, -60, 60, -10, -90, -30]", $false, $false, $false, $false, $false, $true, 1, $false, "[-0.524, -1.048, 1.048, -0.175, -1.571,        -0.524]", 2)
Write-Output "row4 found=$found4"
